$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix header cell B1: was numeric 0, should be text "prompt"
$ws.Range("B1").Value = "prompt"

# Rows 2..101 in column B have a stray leading space in the prompt text.
# Trim it so the text no longer needs xml:space="preserve".
for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = ([string]$val).TrimStart()
    }
}
